$d = $word.ActiveDocument

$pairs = @(
    @('24+19=43', '60-9=51'),
    @('75-38=37', '80+13=93'),
    @('9+78=87', '59-28=31'),
    @('41+18=59', '78-33=45'),
    @('35+2=37', '37+49=86'),
    @('37+12=49', '70+6=76'),
    @('14+33=47', '2+37=39'),
    @('17+9=26', '72-37=35'),
    @('72-33=39', '94-23=71'),
    @('34+46=80', '10-5=5'),
    @('84-69=15', '75-49=26'),
    @('97-9=88', '55-51=4'),
    @('78-30=48', '42-8=34'),
    @('8+75=83', '79-1=78'),
    @('88+5=93', '14+31=45'),
    @('7+81=88', '84-24=60'),
    @('3-0=3', '22+30=52'),
    @('69+4=73', '35+14=49'),
    @('5+62=67', '24+26=50'),
    @('9-4=5', '81-3=78'),
    @('92-23=69', '22+74=96'),
    @('53+7=60', '15+17=32'),
    @('93-59=34', '90-61=29'),
    @('82+7=89', '6+7=13'),
    @('28+61=89', '59+6=65'),
    @('40-9=31', '5+29=34'),
    @('27+45=72', '71-42=29'),
    @('21-12=9', '0+22=22'),
    @('86-53=33', '96-80=16'),
    @('5+3=8', '9+72=81'),
    @('9+47=56', '78+15=93'),
    @('9+66=75', '90-32=58'),
    @('95-94=1', '62+15=77'),
    @('62-21=41', '33+20=53'),
    @('67-25=42', '94-3=91'),
    @('52+2=54', '59-4=55'),
    @('54-37=17', '85-31=54'),
    @('43-18=25', '23-13=10'),
    @('29+23=52', '16+68=84'),
    @('12+76=88', '26-6=20'),
    @('2+21=23', '3+62=65'),
    @('83-4=79', '57-34=23'),
    @('2+8=10', '66-58=8'),
    @('84-27=57', '69+11=80'),
    @('85-8=77', '57+37=94'),
    @('44-26=18', '50+1=51'),
    @('24+3=27', '3+54=57'),
    @('27+0=27', '9+29=38'),
    @('81-57=24', '65-33=32'),
    @('41+6=47', '81-10=71'),
    @('86-56=30', '77-58=19'),
    @('21-5=16', '62-20=42'),
    @('47+12=59', '45-9=36'),
    @('74+22=96', '71-10=61'),
    @('36+32=68', '0+76=76'),
    @('33+65=98', '3+49=52'),
    @('92-63=29', '11+20=31'),
    @('57-31=26', '6+85=91'),
    @('33+31=64', '9+43=52'),
    @('98-62=36', '96-3=93'),
    @('83-0=83', '97-42=55'),
    @('13-2=11', '21+18=39'),
    @('7+8=15', '54+12=66'),
    @('65-32=33', '33+41=74'),
    @('19-10=9', '21+71=92'),
    @('65-11=54', '0+77=77'),
    @('89-78=11', '36+48=84'),
    @('67-10=57', '11+5=16'),
    @('29+14=43', '5+20=25'),
    @('58+2=60', '1+53=54'),
    @('22-10=12', '90-52=38'),
    @('59+39=98', '9+67=76'),
    @('65-20=45', '69-61=8'),
    @('5+47=52', '36+46=82'),
    @('75+10=85', '63-62=1'),
    @('35-16=19', '3+49=52'),
    @('58-11=47', '5+63=68'),
    @('89-20=69', '29-13=16'),
    @('2+84=86', '20+11=31'),
    @('12+0=12', '89-24=65'),
    @('96-57=39', '6+19=25'),
    @('68-13=55', '92-5=87'),
    @('29+41=70', '18-10=8'),
    @('37+54=91', '24+35=59'),
    @('89+6=95', '96-23=73'),
    @('25+47=72', '13+43=56'),
    @('62+23=85', '5+46=51'),
    @('14+40=54', '51+31=82'),
    @('14+45=59', '10+28=38'),
    @('36-23=13', '17+78=95'),
    @('87-12=75', '0+20=20'),
    @('91-76=15', '74+10=84'),
    @('84-48=36', '58+12=70'),
    @('1+40=41', '87-84=3'),
    @('92-31=61', '35+48=83'),
    @('49+34=83', '31+40=71'),
    @('92-87=5', '19-8=11'),
    @('59+29=88', '4+38=42'),
    @('63-19=44', '70-53=17'),
    @('36+41=77', '28+12=40'),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done replacing $($pairs.Count) equations"
